$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 6802
$ws.Range("I21").Value = 8005
$ws.Range("K21").Value = 8005
$ws.Range("M21").Value = -7537

$ws.Range("H23").Value = 6802
$ws.Range("I23").Value = 8005
$ws.Range("K23").Value = 8005
$ws.Range("M23").Value = -7771

$ws.Range("H32").Value = 3844.5
$ws.Range("I32").Value = 3445.5
$ws.Range("J32").Value = 3958.5
$ws.Range("K32").Value = 3445.5
$ws.Range("L32").Value = 3958.5
$ws.Range("M32").Value = -3119.5
$ws.Range("N32").Value = -4610.5

$ws.Range("H40").Value = 7666.3335
$ws.Range("J40").Value = 8999.5
$ws.Range("L40").Value = 8999.5
$ws.Range("N40").Value = -9349.5

$ws.Range("H43").Value = 2683.6191
$ws.Range("I43").Value = 1237.909
$ws.Range("J43").Value = 4273.9
$ws.Range("K43").Value = 1237.909
$ws.Range("L43").Value = 4273.9
$ws.Range("M43").Value = -1168.909
$ws.Range("N43").Value = -4411.9

$ws.Range("H51").Value = 7186.2856
$ws.Range("J51").Value = 7186.2856
$ws.Range("L51").Value = 7186.2856
$ws.Range("N51").Value = -8154.2856

$ws.Range("H64").Value = 8100
$ws.Range("I64").Value = 5250
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 5250
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = -5002
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 8100
$ws.Range("I67").Value = 5250
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 5250
$ws.Range("L67").Value = 10000
$ws.Range("M67").Value = -4392
$ws.Range("N67").Value = -11716

$ws.Range("H112").Value = 4312.3784
$ws.Range("I112").Value = 899.6667
$ws.Range("K112").Value = 2699.0001
$ws.Range("M112").Value = -1591.0001

$ws.Range("H132").Value = 3512.7334
$ws.Range("I132").Value = 3512.7334
$ws.Range("K132").Value = 10538.2002
$ws.Range("M132").Value = -8008.200199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10274.4375
$ws.Range("I32").Value = 5780.256
$ws.Range("K32").Value = 5780.256
$ws.Range("M32").Value = -5493.256

$ws.Range("H97").Value = 4868.5713
$ws.Range("I97").Value = 6329.0713
$ws.Range("J97").Value = 1947.5714
$ws.Range("K97").Value = 6329.0713
$ws.Range("L97").Value = 1947.5714
$ws.Range("M97").Value = -5833.0713
$ws.Range("N97").Value = -2939.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3315.0667
$ws.Range("I107").Value = 2581.4
$ws.Range("K107").Value = 2581.4
$ws.Range("M107").Value = -661.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 167422.22
$ws.Range("J141").Value = 167422.22
$ws.Range("L141").Value = 167422.22
$ws.Range("N141").Value = -177782.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1641.4286
$ws.Range("I50").Value = 1372.5
$ws.Range("K50").Value = 4117.5
$ws.Range("M50").Value = -3636.5

$ws.Range("H53").Value = 1641.4286
$ws.Range("I53").Value = 1372.5
$ws.Range("K53").Value = 4117.5
$ws.Range("M53").Value = -3636.5

$ws.Range("H131").Value = 9472942
$ws.Range("J131").Value = 11115043
$ws.Range("L131").Value = 33345129
$ws.Range("N131").Value = -33355209

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 22999.5
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 22999.5
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 22999.5
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = -23585.5

$ws.Range("H70").Value = 7062.75
$ws.Range("I70").Value = 6198.6
$ws.Range("K70").Value = 6198.6
$ws.Range("M70").Value = -5928.6

$ws.Range("H73").Value = 7062.75
$ws.Range("I73").Value = 6198.6
$ws.Range("K73").Value = 6198.6
$ws.Range("M73").Value = -5262.6

$ws.Range("H126").Value = 4351.346
$ws.Range("I126").Value = 2816.6667
$ws.Range("J126").Value = 4811.75
$ws.Range("K126").Value = 8450.000100000001
$ws.Range("L126").Value = 14435.25
$ws.Range("M126").Value = -5980.000100000001
$ws.Range("N126").Value = -19375.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4730.3335
$ws.Range("I7").Value = 3806.111
$ws.Range("J7").Value = 6116.6665
$ws.Range("K7").Value = 3806.111
$ws.Range("L7").Value = 6116.6665
$ws.Range("M7").Value = -3694.111
$ws.Range("N7").Value = -6340.6665

$ws.Range("H46").Value = 6666.4165
$ws.Range("I46").Value = 2999.6667
$ws.Range("J46").Value = 7888.6665
$ws.Range("K46").Value = 2999.6667
$ws.Range("L46").Value = 7888.6665
$ws.Range("M46").Value = -2811.6667
$ws.Range("N46").Value = -8264.666499999999

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""

$ws.Range("H55").Value = 3037
$ws.Range("I55").Value = 3283.8
$ws.Range("K55").Value = 3283.8
$ws.Range("M55").Value = -3110.8

$ws.Range("H126").Value = 4730.3335
$ws.Range("I126").Value = 3806.111
$ws.Range("J126").Value = 6116.6665
$ws.Range("K126").Value = 11418.333
$ws.Range("L126").Value = 18349.9995
$ws.Range("M126").Value = -8948.332999999999
$ws.Range("N126").Value = -23289.9995

$ws.Range("H132").Value = 10235.096
$ws.Range("I132").Value = 10608.444
$ws.Range("J132").Value = 7995
$ws.Range("K132").Value = 31825.332
$ws.Range("L132").Value = 23985
$ws.Range("M132").Value = -29295.332
$ws.Range("N132").Value = -29045

$ws.Range("H138").Value = 84166
$ws.Range("J138").Value = 84166
$ws.Range("L138").Value = 84166
$ws.Range("N138").Value = -94446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 28999.5
$ws.Range("J37").Value = 28999.5
$ws.Range("L37").Value = 28999.5
$ws.Range("N37").Value = -29405.5

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""

$ws.Range("H81").Value = 10270.818
$ws.Range("I81").Value = 9833
$ws.Range("J81").Value = 10435
$ws.Range("K81").Value = 19666
$ws.Range("L81").Value = 20870
$ws.Range("M81").Value = -18605
$ws.Range("N81").Value = -22992

$ws.Range("H84").Value = 10270.818
$ws.Range("I84").Value = 9833
$ws.Range("J84").Value = 10435
$ws.Range("K84").Value = 98330
$ws.Range("L84").Value = 104350
$ws.Range("M84").Value = -93026
$ws.Range("N84").Value = -114958

$ws.Range("H96").Value = 3386.25
$ws.Range("I96").Value = 3204.0557
$ws.Range("K96").Value = 3204.0557
$ws.Range("M96").Value = -1831.0557

$ws.Range("H113").Value = 936.62067
$ws.Range("I113").Value = 552
$ws.Range("J113").Value = 1109.7
$ws.Range("K113").Value = 1656
$ws.Range("L113").Value = 3329.1
$ws.Range("M113").Value = 514
$ws.Range("N113").Value = -7669.1
